{"js": "const body = context.document.body;\n\n// 1) Paragraph 1: remove \", pour am\u00e9liorer la vitesse du site\" from the\n//    sentence describing why the framework was chosen.\nconst removeTarget = body.search(\", pour am\u00e9liorer la vitesse du site\", { matchCase: false, matchWholeWord: false });\nremoveTarget.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < removeTarget.items.length; i++) {\n  removeTarget.items[i].insertText(\"\", \"Replace\");\n}\nawait context.sync();\n\n// 2) \"Linting and CI\" paragraph: append the extra sentence about not knowing\n//    what linting/CI really is, right after the existing sentence.\nconst anchor = body.search(\"ESLint a mon choix afin de garder le code propre et bien structur\u00e9.\", { matchCase: false, matchWholeWord: false });\nanchor.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < anchor.items.length; i++) {\n  anchor.items[i].insertText(\n    \" Je savais pas ce que c\u2019est mais, j\u2019ai toujours du mal \u00e0 savoir \u00e0 quoi \u00e7a sert mais je penses que c\u2019est pour \u00e9viter les faute de syntaxe etc \",\n    \"End\"\n  );\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Paragraph 1: remove \", pour am\u00e9liorer la vitesse du site\" from the\n#    sentence describing why the framework was chosen.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \", pour am\u00e9liorer la vitesse du site\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2)\n\n# 2) \"Linting and CI\" paragraph: append the extra sentence about not knowing\n#    what linting/CI really is, right after the existing sentence.\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \"ESLint a mon choix afin de garder le code propre et bien structur\u00e9.\"\n$found = $range.Find.Execute()\nif ($found) {\n    # Collapse the found range to its end point (set Start = End) so the\n    # insertion lands right after the matched sentence, regardless of how\n    # the Collapse() direction constants are mapped in this host.\n    $range.Start = $range.End\n    $range.InsertAfter(\" Je savais pas ce que c\u2019est mais, j\u2019ai toujours du mal \u00e0 savoir \u00e0 quoi \u00e7a sert mais je penses que c\u2019est pour \u00e9viter les faute de syntaxe etc \")\n}\n"}
